$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #1")

# Helper: write a literal string into a cell while preserving the cell's
# existing style/number-format. A direct Value assignment on these cells
# would get auto-parsed into a real date serial by the engine for some of
# the "dd-mm-yyyy" strings below, since they are ambiguous with
# mm-dd-yyyy (e.g. "02-02-2017", "06-02-2017", "07-02-2017"). Routing the
# text through a formula cell and then doing a values-only paste avoids
# the date auto-detection while keeping the destination cell's original
# style untouched.
function Set-LiteralText {
    param($sheet, $addr, $text)
    $helper = $sheet.Range("ZZ1")
    $escaped = $text -replace '"', '""'
    $helper.Formula = '="' + $escaped + '"'
    $helper.Copy()
    $sheet.Range($addr).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $helper.Clear()
}

# --- Journal entries (rows 14,15,17-20 first: dates, then descriptions) ---
Set-LiteralText $ws "A14" "30-01-2017"
Set-LiteralText $ws "A15" "31-01-2017"
Set-LiteralText $ws "A17" "06-02-2017"
Set-LiteralText $ws "A18" "07-02-2017"
Set-LiteralText $ws "A19" "13-02-2017"
Set-LiteralText $ws "A20" "14-02-2017"

$ws.Range("B14").Value = "Recherche sur les fragmented activity, plus tutoriel"
$ws.Range("B15").Value = "tentative de faire fonctionner les fragmented activity"
$ws.Range("B17").Value = "utiliser mon cellulaire pour faire foncitonner mon app "
$ws.Range("B18").Value = "suppresion des fragmented activity et utilisation des gridview "
$ws.Range("B19").Value = "erreur avec les gridview commencement avec les layout"
$ws.Range("B20").Value = "travaille sur l'interface  et configuration des bouton pour activé le  Clicke"

$ws.Range("C14").Value = 3
$ws.Range("C15").Value = 2
$ws.Range("C17").Value = 3
$ws.Range("C18").Value = 2
$ws.Range("C19").Value = 3
$ws.Range("C20").Value = 2

# Row 16 was filled in afterwards (it slots chronologically between the
# 31-01 and 06-02 entries above)
Set-LiteralText $ws "A16" "02-02-2017"
$ws.Range("B16").Value = "modification des fragments et essaie d'autres tuto"
$ws.Range("C16").Value = 2

# Row 21 added last
Set-LiteralText $ws "A21" "18-02-2017"
$ws.Range("B21").Value = "recherche pour faire un timer et modification du theme"
$ws.Range("C21").Value = 1

# --- Self-assessment note + grade ---
$ws.Range("B40").Value = 8
$ws.Range("B42").Value = "J'ai mal gerer mon temps mais je suis tout de même satisfait, plusieurs tentative pour faire l'interface qui n'on pas fonctionné. Parcontre je n'ai pas respecter mon itinéraire, j'ai noter que je devais faire l'interface a la fin mais je n'avais pas réfléchit au fait qu'il me faut des objet sur lequel progguer."

# --- Window/selection state: "Iteration #1" becomes the active/selected tab ---
$ws.Range("B42:B47").Select()
$ws.Activate()
